$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 68.443746
$ws.Cells.Item(2, 8).Value = 205.331238
$ws.Cells.Item(2, 9).Value = 0.1596169534001499
$ws.Cells.Item(2, 10).Value = 0.1596169534001499
$ws.Cells.Item(2, 13).Value = 9.423852333333334
$ws.Cells.Item(2, 14).Value = 28.271557
$ws.Cells.Item(2, 15).Value = 0.06654336290212845
$ws.Cells.Item(2, 16).Value = 0.06654336290212845
$ws.Cells.Item(2, 17).Value = 645.0037554441741
$ws.Cells.Item(2, 18).Value = 5805.033798997567
$ws.Cells.Item(2, 19).Value = 0.0106214488554383
$ws.Cells.Item(2, 20).Value = 0.0106214488554383

# Row 3
$ws.Cells.Item(3, 7).Value = 68.443746
$ws.Cells.Item(3, 8).Value = 205.331238
$ws.Cells.Item(3, 9).Value = 0.1596169534001499
$ws.Cells.Item(3, 10).Value = 0.1596169534001499
$ws.Cells.Item(3, 15).Value = 0.3572423751649123
$ws.Cells.Item(3, 16).Value = 0.3572423751649123
$ws.Cells.Item(3, 17).Value = 3462.744645533908
$ws.Cells.Item(3, 18).Value = 31164.70180980517
$ws.Cells.Item(3, 19).Value = 0.05702193954925668
$ws.Cells.Item(3, 20).Value = 0.05702193954925668

# Row 4
$ws.Cells.Item(4, 7).Value = 68.443746
$ws.Cells.Item(4, 8).Value = 205.331238
$ws.Cells.Item(4, 9).Value = 0.1596169534001499
$ws.Cells.Item(4, 10).Value = 0.1596169534001499
$ws.Cells.Item(4, 13).Value = 26.84076266666667
$ws.Cells.Item(4, 14).Value = 80.522288
$ws.Cells.Item(4, 15).Value = 0.1895270158659356
$ws.Cells.Item(4, 16).Value = 0.1895270158659356
$ws.Cells.Item(4, 17).Value = 1837.082342403616
$ws.Cells.Item(4, 18).Value = 16533.74108163255
$ws.Cells.Item(4, 19).Value = 0.03025172485954251
$ws.Cells.Item(4, 20).Value = 0.03025172485954252

# Row 5
$ws.Cells.Item(5, 7).Value = 68.443746
$ws.Cells.Item(5, 8).Value = 205.331238
$ws.Cells.Item(5, 9).Value = 0.1596169534001499
$ws.Cells.Item(5, 10).Value = 0.1596169534001499
$ws.Cells.Item(5, 13).Value = 54.762539
$ws.Cells.Item(5, 14).Value = 164.287617
$ws.Cells.Item(5, 15).Value = 0.3866872460670236
$ws.Cells.Item(5, 16).Value = 0.3866872460670236
$ws.Cells.Item(5, 17).Value = 3748.153309631094
$ws.Cells.Item(5, 18).Value = 33733.37978667984
$ws.Cells.Item(5, 19).Value = 0.06172184013591242
$ws.Cells.Item(5, 20).Value = 0.06172184013591242

# Row 6
$ws.Cells.Item(6, 9).Value = 0.4159650732941736
$ws.Cells.Item(6, 10).Value = 0.4159650732941736
$ws.Cells.Item(6, 13).Value = 9.423852333333334
$ws.Cells.Item(6, 14).Value = 28.271557
$ws.Cells.Item(6, 15).Value = 0.06654336290212845
$ws.Cells.Item(6, 16).Value = 0.06654336290212845
$ws.Cells.Item(6, 17).Value = 1680.893092450799
$ws.Cells.Item(6, 18).Value = 15128.03783205719
$ws.Cells.Item(6, 19).Value = 0.02767971482682465
$ws.Cells.Item(6, 20).Value = 0.02767971482682465

# Row 7
$ws.Cells.Item(7, 9).Value = 0.4159650732941736
$ws.Cells.Item(7, 10).Value = 0.4159650732941736
$ws.Cells.Item(7, 15).Value = 0.3572423751649123
$ws.Cells.Item(7, 16).Value = 0.3572423751649123
$ws.Cells.Item(7, 19).Value = 0.1486003507692574
$ws.Cells.Item(7, 20).Value = 0.1486003507692574

# Row 8
$ws.Cells.Item(8, 9).Value = 0.4159650732941736
$ws.Cells.Item(8, 10).Value = 0.4159650732941736
$ws.Cells.Item(8, 13).Value = 26.84076266666667
$ws.Cells.Item(8, 14).Value = 80.522288
$ws.Cells.Item(8, 15).Value = 0.1895270158659356
$ws.Cells.Item(8, 16).Value = 0.1895270158659356
$ws.Cells.Item(8, 17).Value = 4787.474481420811
$ws.Cells.Item(8, 18).Value = 43087.2703327873
$ws.Cells.Item(8, 19).Value = 0.07883661904589991
$ws.Cells.Item(8, 20).Value = 0.07883661904589991

# Row 9
$ws.Cells.Item(9, 9).Value = 0.4159650732941736
$ws.Cells.Item(9, 10).Value = 0.4159650732941736
$ws.Cells.Item(9, 13).Value = 54.762539
$ws.Cells.Item(9, 14).Value = 164.287617
$ws.Cells.Item(9, 15).Value = 0.3866872460670236
$ws.Cells.Item(9, 16).Value = 0.3866872460670236
$ws.Cells.Item(9, 17).Value = 9767.764845441745
$ws.Cells.Item(9, 18).Value = 87909.88360897571
$ws.Cells.Item(9, 19).Value = 0.1608483886521916
$ws.Cells.Item(9, 20).Value = 0.1608483886521916

# Row 10
$ws.Cells.Item(10, 7).Value = 88.88346833333333
$ws.Cells.Item(10, 8).Value = 266.650405
$ws.Cells.Item(10, 9).Value = 0.2072842188241036
$ws.Cells.Item(10, 10).Value = 0.2072842188241036
$ws.Cells.Item(10, 13).Value = 9.423852333333334
$ws.Cells.Item(10, 14).Value = 28.271557
$ws.Cells.Item(10, 15).Value = 0.06654336290212845
$ws.Cells.Item(10, 16).Value = 0.06654336290212845
$ws.Cells.Item(10, 17).Value = 837.6246804478428
$ws.Cells.Item(10, 18).Value = 7538.622124030585
$ws.Cells.Item(10, 19).Value = 0.01379338899709653
$ws.Cells.Item(10, 20).Value = 0.01379338899709653

# Row 11
$ws.Cells.Item(11, 7).Value = 88.88346833333333
$ws.Cells.Item(11, 8).Value = 266.650405
$ws.Cells.Item(11, 9).Value = 0.2072842188241036
$ws.Cells.Item(11, 10).Value = 0.2072842188241036
$ws.Cells.Item(11, 15).Value = 0.3572423751649123
$ws.Cells.Item(11, 16).Value = 0.3572423751649123
$ws.Cells.Item(11, 17).Value = 4496.842619451785
$ws.Cells.Item(11, 18).Value = 40471.58357506606
$ws.Cells.Item(11, 19).Value = 0.07405070666692619
$ws.Cells.Item(11, 20).Value = 0.07405070666692619

# Row 12
$ws.Cells.Item(12, 7).Value = 88.88346833333333
$ws.Cells.Item(12, 8).Value = 266.650405
$ws.Cells.Item(12, 9).Value = 0.2072842188241036
$ws.Cells.Item(12, 10).Value = 0.2072842188241036
$ws.Cells.Item(12, 13).Value = 26.84076266666667
$ws.Cells.Item(12, 14).Value = 80.522288
$ws.Cells.Item(12, 15).Value = 0.1895270158659356
$ws.Cells.Item(12, 16).Value = 0.1895270158659356
$ws.Cells.Item(12, 17).Value = 2385.700078525182
$ws.Cells.Item(12, 18).Value = 21471.30070672664
$ws.Cells.Item(12, 19).Value = 0.03928595942983394
$ws.Cells.Item(12, 20).Value = 0.03928595942983395

# Row 13
$ws.Cells.Item(13, 7).Value = 88.88346833333333
$ws.Cells.Item(13, 8).Value = 266.650405
$ws.Cells.Item(13, 9).Value = 0.2072842188241036
$ws.Cells.Item(13, 10).Value = 0.2072842188241036
$ws.Cells.Item(13, 13).Value = 54.762539
$ws.Cells.Item(13, 14).Value = 164.287617
$ws.Cells.Item(13, 15).Value = 0.3866872460670236
$ws.Cells.Item(13, 16).Value = 0.3866872460670236
$ws.Cells.Item(13, 17).Value = 4867.484401059431
$ws.Cells.Item(13, 18).Value = 43807.35960953488
$ws.Cells.Item(13, 19).Value = 0.08015416373024692
$ws.Cells.Item(13, 20).Value = 0.08015416373024692

# Row 14
$ws.Cells.Item(14, 7).Value = 93.106949
$ws.Cells.Item(14, 8).Value = 279.320847
$ws.Cells.Item(14, 9).Value = 0.2171337544815728
$ws.Cells.Item(14, 10).Value = 0.2171337544815728
$ws.Cells.Item(14, 13).Value = 9.423852333333334
$ws.Cells.Item(14, 14).Value = 28.271557
$ws.Cells.Item(14, 15).Value = 0.06654336290212845
$ws.Cells.Item(14, 16).Value = 0.06654336290212845
$ws.Cells.Item(14, 17).Value = 877.4261385831977
$ws.Cells.Item(14, 18).Value = 7896.83524724878
$ws.Cells.Item(14, 19).Value = 0.01444881022276896
$ws.Cells.Item(14, 20).Value = 0.01444881022276896

# Row 15
$ws.Cells.Item(15, 7).Value = 93.106949
$ws.Cells.Item(15, 8).Value = 279.320847
$ws.Cells.Item(15, 9).Value = 0.2171337544815728
$ws.Cells.Item(15, 10).Value = 0.2171337544815728
$ws.Cells.Item(15, 15).Value = 0.3572423751649123
$ws.Cells.Item(15, 16).Value = 0.3572423751649123
$ws.Cells.Item(15, 17).Value = 4710.519338198535
$ws.Cells.Item(15, 18).Value = 42394.67404378682
$ws.Cells.Item(15, 19).Value = 0.077569378179472
$ws.Cells.Item(15, 20).Value = 0.077569378179472

# Row 16
$ws.Cells.Item(16, 7).Value = 93.106949
$ws.Cells.Item(16, 8).Value = 279.320847
$ws.Cells.Item(16, 9).Value = 0.2171337544815728
$ws.Cells.Item(16, 10).Value = 0.2171337544815728
$ws.Cells.Item(16, 13).Value = 26.84076266666667
$ws.Cells.Item(16, 14).Value = 80.522288
$ws.Cells.Item(16, 15).Value = 0.1895270158659356
$ws.Cells.Item(16, 16).Value = 0.1895270158659356
$ws.Cells.Item(16, 17).Value = 2499.061520726437
$ws.Cells.Item(16, 18).Value = 22491.55368653794
$ws.Cells.Item(16, 19).Value = 0.04115271253065921
$ws.Cells.Item(16, 20).Value = 0.04115271253065922

# Row 17
$ws.Cells.Item(17, 7).Value = 93.106949
$ws.Cells.Item(17, 8).Value = 279.320847
$ws.Cells.Item(17, 9).Value = 0.2171337544815728
$ws.Cells.Item(17, 10).Value = 0.2171337544815728
$ws.Cells.Item(17, 13).Value = 54.762539
$ws.Cells.Item(17, 14).Value = 164.287617
$ws.Cells.Item(17, 15).Value = 0.3866872460670236
$ws.Cells.Item(17, 16).Value = 0.3866872460670236
$ws.Cells.Item(17, 17).Value = 5098.772925783511
$ws.Cells.Item(17, 18).Value = 45888.9563320516
$ws.Cells.Item(17, 19).Value = 0.08396285354867265
$ws.Cells.Item(17, 20).Value = 0.08396285354867265
